$d = $word.ActiveDocument

$d.Content.Find.Execute("75-39=", $true, $false, $false, $false, $false, $true, 1, $false, "40-2=", 2) | Out-Null
$d.Content.Find.Execute("68-5=", $true, $false, $false, $false, $false, $true, 1, $false, "4+85=", 2) | Out-Null
$d.Content.Find.Execute("63-18=", $true, $false, $false, $false, $false, $true, 1, $false, "69-46=", 2) | Out-Null
$d.Content.Find.Execute("76-5=", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=", 2) | Out-Null
$d.Content.Find.Execute("40+11=", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=", 2) | Out-Null
$d.Content.Find.Execute("78-9=", $true, $false, $false, $false, $false, $true, 1, $false, "57+41=", 2) | Out-Null
$d.Content.Find.Execute("31+59=", $true, $false, $false, $false, $false, $true, 1, $false, "2+81=", 2) | Out-Null
$d.Content.Find.Execute("35-24=", $true, $false, $false, $false, $false, $true, 1, $false, "67+18=", 2) | Out-Null
$d.Content.Find.Execute("94-67=", $true, $false, $false, $false, $false, $true, 1, $false, "94+0=", 2) | Out-Null
$d.Content.Find.Execute("64-23=", $true, $false, $false, $false, $false, $true, 1, $false, "30+17=", 2) | Out-Null
$d.Content.Find.Execute("36+1=", $true, $false, $false, $false, $false, $true, 1, $false, "39+55=", 2) | Out-Null
$d.Content.Find.Execute("51-5=", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=", 2) | Out-Null
$d.Content.Find.Execute("1+6=", $true, $false, $false, $false, $false, $true, 1, $false, "15+67=", 2) | Out-Null
$d.Content.Find.Execute("8+70=", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("23+23=", $true, $false, $false, $false, $false, $true, 1, $false, "68+18=", 2) | Out-Null
$d.Content.Find.Execute("41+22=", $true, $false, $false, $false, $false, $true, 1, $false, "42-24=", 2) | Out-Null
$d.Content.Find.Execute("17-6=", $true, $false, $false, $false, $false, $true, 1, $false, "12+25=", 2) | Out-Null
$d.Content.Find.Execute("46-38=", $true, $false, $false, $false, $false, $true, 1, $false, "69+5=", 2) | Out-Null
$d.Content.Find.Execute("3+49=", $true, $false, $false, $false, $false, $true, 1, $false, "49+1=", 2) | Out-Null
$d.Content.Find.Execute("65-61=", $true, $false, $false, $false, $false, $true, 1, $false, "1+60=", 2) | Out-Null
$d.Content.Find.Execute("97-13=", $true, $false, $false, $false, $false, $true, 1, $false, "69-24=", 2) | Out-Null
$d.Content.Find.Execute("33+5=", $true, $false, $false, $false, $false, $true, 1, $false, "42+24=", 2) | Out-Null
$d.Content.Find.Execute("88-46=", $true, $false, $false, $false, $false, $true, 1, $false, "70+15=", 2) | Out-Null
$d.Content.Find.Execute("28+12=", $true, $false, $false, $false, $false, $true, 1, $false, "95-29=", 2) | Out-Null
$d.Content.Find.Execute("47-47=", $true, $false, $false, $false, $false, $true, 1, $false, "48+19=", 2) | Out-Null
$d.Content.Find.Execute("14+57=", $true, $false, $false, $false, $false, $true, 1, $false, "30-11=", 2) | Out-Null
$d.Content.Find.Execute("19+41=", $true, $false, $false, $false, $false, $true, 1, $false, "81-22=", 2) | Out-Null
$d.Content.Find.Execute("61-10=", $true, $false, $false, $false, $false, $true, 1, $false, "19+80=", 2) | Out-Null
$d.Content.Find.Execute("61-1=", $true, $false, $false, $false, $false, $true, 1, $false, "17-2=", 2) | Out-Null
$d.Content.Find.Execute("87-8=", $true, $false, $false, $false, $false, $true, 1, $false, "7+39=", 2) | Out-Null
$d.Content.Find.Execute("11+81=", $true, $false, $false, $false, $false, $true, 1, $false, "45+44=", 2) | Out-Null
$d.Content.Find.Execute("15+64=", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=", 2) | Out-Null
$d.Content.Find.Execute("20+31=", $true, $false, $false, $false, $false, $true, 1, $false, "51+3=", 2) | Out-Null
$d.Content.Find.Execute("93-84=", $true, $false, $false, $false, $false, $true, 1, $false, "82-65=", 2) | Out-Null
$d.Content.Find.Execute("23-0=", $true, $false, $false, $false, $false, $true, 1, $false, "72-1=", 2) | Out-Null
$d.Content.Find.Execute("57-56=", $true, $false, $false, $false, $false, $true, 1, $false, "74-41=", 2) | Out-Null
$d.Content.Find.Execute("2+4=", $true, $false, $false, $false, $false, $true, 1, $false, "31+53=", 2) | Out-Null
$d.Content.Find.Execute("15+51=", $true, $false, $false, $false, $false, $true, 1, $false, "40-11=", 2) | Out-Null
$d.Content.Find.Execute("68-1=", $true, $false, $false, $false, $false, $true, 1, $false, "37-32=", 2) | Out-Null
$d.Content.Find.Execute("99-98=", $true, $false, $false, $false, $false, $true, 1, $false, "62-19=", 2) | Out-Null
$d.Content.Find.Execute("43-28=", $true, $false, $false, $false, $false, $true, 1, $false, "67-49=", 2) | Out-Null
$d.Content.Find.Execute("0+20=", $true, $false, $false, $false, $false, $true, 1, $false, "45+34=", 2) | Out-Null
$d.Content.Find.Execute("65-20=", $true, $false, $false, $false, $false, $true, 1, $false, "83+14=", 2) | Out-Null
$d.Content.Find.Execute("26+57=", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=", 2) | Out-Null
$d.Content.Find.Execute("38+32=", $true, $false, $false, $false, $false, $true, 1, $false, "31+51=", 2) | Out-Null
$d.Content.Find.Execute("38-5=", $true, $false, $false, $false, $false, $true, 1, $false, "42+16=", 2) | Out-Null
$d.Content.Find.Execute("60-11=", $true, $false, $false, $false, $false, $true, 1, $false, "72-61=", 2) | Out-Null
$d.Content.Find.Execute("6+17=", $true, $false, $false, $false, $false, $true, 1, $false, "73-20=", 2) | Out-Null
$d.Content.Find.Execute("4+62=", $true, $false, $false, $false, $false, $true, 1, $false, "17-4=", 2) | Out-Null
$d.Content.Find.Execute("5+85=", $true, $false, $false, $false, $false, $true, 1, $false, "83-60=", 2) | Out-Null
$d.Content.Find.Execute("22+6=", $true, $false, $false, $false, $false, $true, 1, $false, "76-4=", 2) | Out-Null
$d.Content.Find.Execute("23+0=", $true, $false, $false, $false, $false, $true, 1, $false, "20+58=", 2) | Out-Null
$d.Content.Find.Execute("68+21=", $true, $false, $false, $false, $false, $true, 1, $false, "93+6=", 2) | Out-Null
$d.Content.Find.Execute("49-18=", $true, $false, $false, $false, $false, $true, 1, $false, "2+30=", 2) | Out-Null
$d.Content.Find.Execute("65+8=", $true, $false, $false, $false, $false, $true, 1, $false, "66+0=", 2) | Out-Null
$d.Content.Find.Execute("39-4=", $true, $false, $false, $false, $false, $true, 1, $false, "59-25=", 2) | Out-Null
$d.Content.Find.Execute("67-53=", $true, $false, $false, $false, $false, $true, 1, $false, "97-54=", 2) | Out-Null
$d.Content.Find.Execute("84-82=", $true, $false, $false, $false, $false, $true, 1, $false, "87-3=", 2) | Out-Null
$d.Content.Find.Execute("65-25=", $true, $false, $false, $false, $false, $true, 1, $false, "58-16=", 2) | Out-Null
$d.Content.Find.Execute("42+50=", $true, $false, $false, $false, $false, $true, 1, $false, "75+13=", 2) | Out-Null
$d.Content.Find.Execute("29-17=", $true, $false, $false, $false, $false, $true, 1, $false, "85-26=", 2) | Out-Null
$d.Content.Find.Execute("29+5=", $true, $false, $false, $false, $false, $true, 1, $false, "49+5=", 2) | Out-Null
$d.Content.Find.Execute("98-77=", $true, $false, $false, $false, $false, $true, 1, $false, "11+84=", 2) | Out-Null
$d.Content.Find.Execute("29+46=", $true, $false, $false, $false, $false, $true, 1, $false, "15+15=", 2) | Out-Null
$d.Content.Find.Execute("69-45=", $true, $false, $false, $false, $false, $true, 1, $false, "15+12=", 2) | Out-Null
$d.Content.Find.Execute("88-76=", $true, $false, $false, $false, $false, $true, 1, $false, "52+41=", 2) | Out-Null
$d.Content.Find.Execute("28+34=", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=", 2) | Out-Null
$d.Content.Find.Execute("37+48=", $true, $false, $false, $false, $false, $true, 1, $false, "84-74=", 2) | Out-Null
$d.Content.Find.Execute("32-25=", $true, $false, $false, $false, $false, $true, 1, $false, "16-12=", 2) | Out-Null
$d.Content.Find.Execute("94-57=", $true, $false, $false, $false, $false, $true, 1, $false, "50+35=", 2) | Out-Null
$d.Content.Find.Execute("22-15=", $true, $false, $false, $false, $false, $true, 1, $false, "99-91=", 2) | Out-Null
$d.Content.Find.Execute("24+29=", $true, $false, $false, $false, $false, $true, 1, $false, "63+31=", 2) | Out-Null
$d.Content.Find.Execute("64-37=", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=", 2) | Out-Null
$d.Content.Find.Execute("86-53=", $true, $false, $false, $false, $false, $true, 1, $false, "3+38=", 2) | Out-Null
$d.Content.Find.Execute("69-49=", $true, $false, $false, $false, $false, $true, 1, $false, "22+75=", 2) | Out-Null
$d.Content.Find.Execute("99-20=", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=", 2) | Out-Null
$d.Content.Find.Execute("8+2=", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("73-29=", $true, $false, $false, $false, $false, $true, 1, $false, "98-11=", 2) | Out-Null
$d.Content.Find.Execute("98-88=", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=", 2) | Out-Null
$d.Content.Find.Execute("88-50=", $true, $false, $false, $false, $false, $true, 1, $false, "75-21=", 2) | Out-Null
$d.Content.Find.Execute("3+78=", $true, $false, $false, $false, $false, $true, 1, $false, "96-11=", 2) | Out-Null
$d.Content.Find.Execute("13+25=", $true, $false, $false, $false, $false, $true, 1, $false, "78-14=", 2) | Out-Null
$d.Content.Find.Execute("30-5=", $true, $false, $false, $false, $false, $true, 1, $false, "29+31=", 2) | Out-Null
$d.Content.Find.Execute("21+44=", $true, $false, $false, $false, $false, $true, 1, $false, "50+20=", 2) | Out-Null
$d.Content.Find.Execute("80-16=", $true, $false, $false, $false, $false, $true, 1, $false, "5+87=", 2) | Out-Null
$d.Content.Find.Execute("87-77=", $true, $false, $false, $false, $false, $true, 1, $false, "98-93=", 2) | Out-Null
$d.Content.Find.Execute("6+65=", $true, $false, $false, $false, $false, $true, 1, $false, "59-53=", 2) | Out-Null
$d.Content.Find.Execute("24+58=", $true, $false, $false, $false, $false, $true, 1, $false, "29+19=", 2) | Out-Null
$d.Content.Find.Execute("18+37=", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=", 2) | Out-Null
$d.Content.Find.Execute("38+27=", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=", 2) | Out-Null
$d.Content.Find.Execute("52+24=", $true, $false, $false, $false, $false, $true, 1, $false, "39+13=", 2) | Out-Null
$d.Content.Find.Execute("97-61=", $true, $false, $false, $false, $false, $true, 1, $false, "73+23=", 2) | Out-Null
$d.Content.Find.Execute("56+8=", $true, $false, $false, $false, $false, $true, 1, $false, "72+27=", 2) | Out-Null
$d.Content.Find.Execute("67-14=", $true, $false, $false, $false, $false, $true, 1, $false, "22+29=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $false, $false, $false, $false, $true, 1, $false, "49-37=", 2) | Out-Null
$d.Content.Find.Execute("34+43=", $true, $false, $false, $false, $false, $true, 1, $false, "23-17=", 2) | Out-Null
$d.Content.Find.Execute("86-0=", $true, $false, $false, $false, $false, $true, 1, $false, "87+4=", 2) | Out-Null
$d.Content.Find.Execute("83-56=", $true, $false, $false, $false, $false, $true, 1, $false, "2+1=", 2) | Out-Null
$d.Content.Find.Execute("55+3=", $true, $false, $false, $false, $false, $true, 1, $false, "6+32=", 2) | Out-Null
$d.Content.Find.Execute("58+16=", $true, $false, $false, $false, $false, $true, 1, $false, "30+13=", 2) | Out-Null
